# "Update countries & provincias Spain" - refresh the COVID-19 per-country
# stats on sheet "Pais" with a newer snapshot (22 May 2020, 17:05 instead of
# 16:35). The sheet is kept sorted descending by column B ("Casos totales"),
# so a few rows whose totals changed relative order swap places with their
# neighbour (Nepal/Tanzania, and the Groenlandia/Montserrat/Seychelles trio).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 17:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1625039
$ws.Range("C4").Value = 4137
$ws.Range("D4").Value = 382987
$ws.Range("E4").Value = 1145526
$ws.Range("G4").Value = 172
$ws.Range("H4").Value = 96526

# Alemania (row 11)
$ws.Range("B11").Value = 179347
$ws.Range("C11").Value = 326
$ws.Range("E11").Value = 12027
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 8320

# India (row 14)
$ws.Range("B14").Value = 123711
$ws.Range("C14").Value = 5485
$ws.Range("D14").Value = 50857
$ws.Range("E14").Value = 69178
$ws.Range("G14").Value = 92
$ws.Range("H14").Value = 3676

# Serbia (row 49)
$ws.Range("D49").Value = 5541
$ws.Range("E49").Value = 5246

# Noruega (row 55)
$ws.Range("B55").Value = 8322
$ws.Range("C55").Value = 13
$ws.Range("E55").Value = 360

# Kazajistan (row 57)
$ws.Range("D57").Value = 4096
$ws.Range("E57").Value = 3466

# Sri Lanka (row 106)
$ws.Range("B106").Value = 1060
$ws.Range("C106").Value = 12
$ws.Range("E106").Value = 431

# Nepal moves above Tanzania (new totals re-sort them)
$ws.Range("A130").Value = "Nepal"
$ws.Range("B130").Value = 516
$ws.Range("C130").Value = 59
$ws.Range("D130").Value = 70
$ws.Range("E130").Value = 443
$ws.Range("H130").Value = 3

$ws.Range("A131").Value = "Tanzania"
$ws.Range("B131").Value = 509
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 183
$ws.Range("E131").Value = 305
$ws.Range("H131").Value = 21

# Cabo Verde (row 140)
$ws.Range("B140").Value = 362
$ws.Range("C140").Value = 6
$ws.Range("E140").Value = 264

# Liberia (row 150)
$ws.Range("B150").Value = 249
$ws.Range("C150").Value = 9
$ws.Range("D150").Value = 136
$ws.Range("E150").Value = 89
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 24

# Mozambique (row 157)
$ws.Range("B157").Value = 164
$ws.Range("C157").Value = 2
$ws.Range("E157").Value = 116

# Gibraltar (row 160)
$ws.Range("D160").Value = 147
$ws.Range("E160").Value = 4

# Montserrat / Seychelles / Groenlandia cyclic re-sort (rows 209-211)
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Groenlandia"
